$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 98
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01-09-2021"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 202
$ws.Cells.Item($row, 3).Value = 50
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 152
$ws.Cells.Item($row, 11).Value = 0
